# "womens world cup 2023 data"
# Adds 5 new countries (New Zealand, Philippines, Zambia, Haiti, Vietnam)
# to the localization table on the "Team" sheet, grows the Table6
# ListObject to cover them, and makes "Team" the active/selected sheet
# with the new rows selected (mirroring the saved workbook view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Team")

# New rows: key -> en, es, it, fr, de, nl, ja, fa (columns B..I)
$translations = @{
  81 = @("New Zealand", "Nueva Zelanda", "Nuova Zelanda", "Nouvelle-Zélande", "Neuseeland", "Nieuw-Zeeland", "ニュージーランド", "نیوزلند")
  82 = @("Philippines", "Filipinas", "Filippine", "Philippines", "Philippinen", "Filippijnen", "フィリピン", "وابسته به فیلیپین")
  83 = @("Zambia", "Zambia", "Zambia", "Zambie", "Sambia", "Zambia", "ザンビア", "زامبیا")
  84 = @("Haiti", "Haití", "Haiti", "Haïti", "Haiti", "Haïti", "ハイチ", "هائیتی")
  85 = @("Vietnam", "Vietnam", "Vietnam", "Vietnam", "Vietnam", "Vietnam", "ベトナム", "ویتنام")
}
$keys = @{ 81 = "nzl"; 82 = "phi"; 83 = "zam"; 84 = "hai"; 85 = "vie" }

# Fill in the translation columns (B..I) first, row by row...
foreach ($r in 81..85) {
  $vals = $translations[$r]
  for ($c = 2; $c -le 9; $c++) {
    $ws.Cells.Item($r, $c).Value = $vals[$c - 2]
  }
}

# ...then the key column (A), matching the order the shared strings were
# originally authored in.
foreach ($r in 81..85) {
  $ws.Cells.Item($r, 1).Value = $keys[$r]
}

# Grow the "Table6" listobject so the new rows are part of the table.
$lo = $ws.ListObjects.Item(1)
[void]$lo.Resize($ws.Range("A1:I85"))

# Switch to the Team sheet and select the newly-added rows, matching the
# saved view state (tabSelected/activeTab move from "Fonts" to "Team").
[void]$ws.Activate()
[void]$ws.Range("A81:A85").Select()
